$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking price strings
# (e.g. "1.001", "28.952.06") are stored verbatim as text, matching the
# source data which uses "." as a thousands-style separator, not a decimal point.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.952.06"
$ws.Range("E2").Value = "  -2.65%  "
$ws.Range("D3").Value = "1.884.63"
$ws.Range("E3").Value = "  -3.34%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "329.72"
$ws.Range("E5").Value = "  -3.75%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "0.4585"
$ws.Range("E7").Value = "  -4.36%  "
$ws.Range("D8").Value = "0.4103"
$ws.Range("E8").Value = "  -0.95%  "
$ws.Range("E9").Value = "  -2.04%  "
$ws.Range("D10").Value = "0.07964"
$ws.Range("E10").Value = "  -3.86%  "
$ws.Range("D11").Value = "0.9938"
$ws.Range("E11").Value = "  -5.28%  "
$ws.Range("D12").Value = "21.64"
$ws.Range("E12").Value = "  -4.87%  "
$ws.Range("D13").Value = "1.914.84"
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").Value = "5.910"
$ws.Range("E14").Value = "  -4.35%  "
$ws.Range("D15").Value = "7.066"
$ws.Range("E15").Value = "  -5.11%  "
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "88.40"
$ws.Range("E17").Value = "  -4.89%  "
$ws.Range("D18").Value = "0.06573"
$ws.Range("E18").Value = "  -2.23%  "
$ws.Range("D19").Value = "0.00001027"
$ws.Range("E19").Value = "  -3.78%  "
$ws.Range("D20").Value = "17.41"
$ws.Range("E20").Value = "  -3.64%  "
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "28.923.20"
$ws.Range("E22").Value = "  -2.63%  "
$ws.Range("D23").Value = "5.408"
$ws.Range("E23").Value = "  -3.99%  "
$ws.Range("D24").Value = "11.49"
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").Value = "2.197"
$ws.Range("E25").Value = "  -3.52%  "
$ws.Range("D26").Value = "2.119.13"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("D27").Value = "156.12"
$ws.Range("E27").Value = "  -3.62%  "
$ws.Range("D28").Value = "19.56"
$ws.Range("E28").Value = "  -2.97%  "
$ws.Range("D29").Value = "2.082"
$ws.Range("E29").Value = "  -5.39%  "
$ws.Range("D30").Value = "5.477"
$ws.Range("E30").Value = "  -3.05%  "
$ws.Range("D31").Value = "117.45"
$ws.Range("E31").Value = "  -4.28%  "
$ws.Range("D32").Value = "1.028"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").Value = "0.09320"
$ws.Range("E33").Value = "  -3.51%  "
$ws.Range("D34").Value = "1.402"
$ws.Range("E34").Value = "  -4.96%  "
$ws.Range("E35").Value = "  -4.39%  "
$ws.Range("D36").Value = "5.286"
$ws.Range("E36").Value = "  -3.83%  "
$ws.Range("D37").Value = "0.06052"
$ws.Range("E37").Value = "  -3.12%  "
$ws.Range("D38").Value = "0.02226"
$ws.Range("E38").Value = "  -3.92%  "
$ws.Range("D39").Value = "8.325"
$ws.Range("E39").Value = "  -4.49%  "
$ws.Range("D40").Value = "1.174"
$ws.Range("E40").Value = "  -2.20%  "
$ws.Range("D41").Value = "1.000"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "0.5782"
$ws.Range("E42").Value = "  -5.42%  "
$ws.Range("D43").Value = "0.1823"
$ws.Range("E43").Value = "  -4.52%  "
$ws.Range("D44").Value = "10.06"
$ws.Range("E44").Value = "  -6.54%  "
$ws.Range("D45").Value = "1.244"
$ws.Range("E45").Value = "  -3.00%  "
$ws.Range("D46").Value = "0.07515"
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "2.250"
$ws.Range("E47").Value = "  -3.36%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "0.5450"
$ws.Range("E48").Value = "  -4.88%  "
$ws.Range("D49").Value = "11.95"
$ws.Range("E49").Value = "  -5.60%  "
$ws.Range("D50").Value = "1.899"
$ws.Range("E50").Value = "  -5.07%  "
$ws.Range("D51").Value = "111.38"
$ws.Range("E51").Value = "  -2.49%  "
